# Ran model for 1/19/2021
# Fill in "Beat Vegas?" (column G) results for the existing 1/18/2021 games
# and append the two new games played on 1/19/2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G results for rows 90-96 (games on 1/18/2021)
$ws.Range("G90").Value = "No"
$ws.Range("G91").Value = "Yes"
$ws.Range("G92").Value = "Yes"
$ws.Range("G93").Value = "Yes"
$ws.Range("G94").Value = "No"
$ws.Range("G95").Value = "Yes"
$ws.Range("G96").Value = "No"

# New game rows for 1/19/2021
$ws.Range("A97").Value = 44215
$ws.Range("A97").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B97").Value = "DEN"
$ws.Range("C97").Value = "OKC"
$ws.Range("D97").Value = -9.5
$ws.Range("E97").Value = -1.6
$ws.Range("F97").Value = -7.9

$ws.Range("A98").Value = 44215
$ws.Range("A98").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B98").Value = "UTA"
$ws.Range("C98").Value = "NOP"
$ws.Range("D98").Value = -6
$ws.Range("E98").Value = -5.6
$ws.Range("F98").Value = -0.40000000000000041

# Match the author's final on-screen selection
[void]$ws.Range("I8").Select()
